# Insert a new weekly record row right before the current row 127.
# This shifts the existing rows 127-140 down to 128-141 and extends
# the used range of the sheet to A1:R141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new Pepino ensalada record.
$ws.Range("A127").Value = 7
$ws.Range("B127").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C127").Value = "Ñuble"
$ws.Range("D127").Value = 44449
$ws.Range("E127").Value = 16
$ws.Range("F127").Value = 100112043
$ws.Range("G127").Value = "Pepino ensalada"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 160
$ws.Range("K127").Value = 16000
$ws.Range("L127").Value = 17000
$ws.Range("M127").Value = 16500
$ws.Range("N127").Value = "$/caja 60 unidades"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 275
$ws.Range("Q127").Value = 60
$ws.Range("R127").Value = "Hortaliza"
